# Insert a new data row for Brócoli (Vega Modelo de Temuco) right before the
# current row 484. This shifts all existing rows 484..594 down to 485..595,
# preserving their data and formatting (including the date number format on
# column D), and extends the sheet's used range to A1:R595.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(484).EntireRow.Insert()

$ws.Cells.Item(484, 1).Value = 10
$ws.Cells.Item(484, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(484, 3).Value = "La Araucanía"
$ws.Cells.Item(484, 4).Value = 44964
$ws.Cells.Item(484, 5).Value = 9
$ws.Cells.Item(484, 6).Value = 100112023
$ws.Cells.Item(484, 7).Value = "Brócoli"
$ws.Cells.Item(484, 8).Value = "Sin especificar"
$ws.Cells.Item(484, 9).Value = "Primera"
$ws.Cells.Item(484, 10).Value = 500
$ws.Cells.Item(484, 11).Value = 1200
$ws.Cells.Item(484, 12).Value = 1300
$ws.Cells.Item(484, 13).Value = 1260
$ws.Cells.Item(484, 14).Value = "$/unidad"
$ws.Cells.Item(484, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(484, 16).Value = 1260
$ws.Cells.Item(484, 17).Value = 1
$ws.Cells.Item(484, 18).Value = "Hortaliza"
